# Generate Report for Handback
# The 17ae5c38-1829-4c49-929c-4025ddce299a.md file's handback transform failed
# because the handback file name didn't match the handoff file name.
# Update the Overview sheet's status for both locales, update each locale
# sheet's Status + Error Detail columns, and widen the Error Detail column
# so the message is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: ufonam52.vjk is different with handoff file name: 17ae5c38-1829-4c49-929c-4025ddce299a.c5053e22d369b738abc37f53ce0a09b78f93e5ad.zh-cn."
$deError = "Handback file name: ufonam52.vjk is different with handoff file name: 17ae5c38-1829-4c49-929c-4025ddce299a.c5053e22d369b738abc37f53ce0a09b78f93e5ad.de-de."

# --- Overview sheet: update zh-cn / de-de status for the 17ae5c38 row (row 7) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E7").Value = $newStatus
$wsOverview.Range("F7").Value = $newStatus

# --- zh-cn sheet: Status (C7) + Error Detail (P7), widen Error Detail column ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C7").Value = $newStatus
$wsZh.Range("P7").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: Status (C7) + Error Detail (P7), widen Error Detail column ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C7").Value = $newStatus
$wsDe.Range("P7").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = 39.17
